$wb = $excel.ActiveWorkbook

# ---- Sheet ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H17").Value = 722.2105
$ws.Range("J17").Value = 797.625
$ws.Range("L17").Value = 2392.875
$ws.Range("N17").Value = -2728.875
$ws.Range("H32").Value = 6603.7646
$ws.Range("I32").Value = 7121.2
$ws.Range("J32").Value = 5864.5713
$ws.Range("K32").Value = 7121.2
$ws.Range("L32").Value = 5864.5713
$ws.Range("M32").Value = -6795.2
$ws.Range("N32").Value = -6516.5713
$ws.Range("H57").Value = 22500
$ws.Range("I57").Value = 15000
$ws.Range("J57").Value = 30000
$ws.Range("K57").Value = 45000
$ws.Range("L57").Value = 90000
$ws.Range("M57").Value = -44501
$ws.Range("N57").Value = -90998
$ws.Range("H125").Value = 3563.4119
$ws.Range("J125").Value = 3399.4
$ws.Range("L125").Value = 30594.6
$ws.Range("N125").Value = -35514.60000000001
$ws.Range("H137").Value = 3260.5
$ws.Range("I137").Value = 1997.5555
$ws.Range("K137").Value = 5992.666499999999
$ws.Range("M137").Value = -3442.666499999999

# ---- Sheet ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 9022.388999999999
$ws.Range("J32").Value = 29601
$ws.Range("L32").Value = 29601
$ws.Range("N32").Value = -30175
$ws.Range("H122").Value = 1738.7727
$ws.Range("I122").Value = 1786.8422
$ws.Range("J122").Value = 1434.3334
$ws.Range("K122").Value = 5360.5266
$ws.Range("L122").Value = 4303.0002
$ws.Range("M122").Value = -2910.5266
$ws.Range("N122").Value = -9203.0002

# ---- Sheet BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H87").Value = 66771.664
$ws.Range("I87").Value = 60157.5
$ws.Range("J87").Value = 80000
$ws.Range("K87").Value = 60157.5
$ws.Range("L87").Value = 80000
$ws.Range("M87").Value = -58909.5
$ws.Range("N87").Value = -82496
$ws.Range("H90").Value = 66771.664
$ws.Range("I90").Value = 60157.5
$ws.Range("J90").Value = 80000
$ws.Range("K90").Value = 180472.5
$ws.Range("L90").Value = 240000
$ws.Range("M90").Value = -174232.5
$ws.Range("N90").Value = -252480
$ws.Range("H107").Value = 1629.4667
$ws.Range("I107").Value = 1341.8125
$ws.Range("J107").Value = 2337.5386
$ws.Range("K107").Value = 1341.8125
$ws.Range("L107").Value = 2337.5386
$ws.Range("M107").Value = 578.1875
$ws.Range("N107").Value = -6177.5386

# ---- Sheet CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 2611.9207
$ws.Range("I31").Value = 2733.96
$ws.Range("J31").Value = 2531.6316
$ws.Range("K31").Value = 2733.96
$ws.Range("L31").Value = 2531.6316
$ws.Range("M31").Value = -2438.96
$ws.Range("N31").Value = -3121.6316
$ws.Range("H34").Value = 2611.9207
$ws.Range("I34").Value = 2733.96
$ws.Range("J34").Value = 2531.6316
$ws.Range("K34").Value = 2733.96
$ws.Range("L34").Value = 2531.6316
$ws.Range("M34").Value = -2531.96
$ws.Range("N34").Value = -2935.6316
$ws.Range("H58").Value = 14113.556
$ws.Range("I58").Value = 14113.556
$ws.Range("J58").Value = 0
$ws.Range("K58").Value = 14113.556
$ws.Range("L58").Value = 0
$ws.Range("M58").ClearContents()
$ws.Range("N58").Value = -13910.556
$ws.Range("H62").Value = 5577.5884
$ws.Range("I62").Value = 4946.357
$ws.Range("K62").Value = 4946.357
$ws.Range("M62").Value = -4322.357
$ws.Range("H65").Value = 5577.5884
$ws.Range("I65").Value = 4946.357
$ws.Range("K65").Value = 24731.785
$ws.Range("M65").Value = -21611.785
$ws.Range("H74").Value = 73499.336
$ws.Range("J74").Value = 99999.664
$ws.Range("L74").Value = 99999.664
$ws.Range("N74").Value = -101747.664
$ws.Range("H77").Value = 73499.336
$ws.Range("J77").Value = 99999.664
$ws.Range("L77").Value = 299998.992
$ws.Range("N77").Value = -308734.992
$ws.Range("H134").Value = 3329.7917
$ws.Range("I134").Value = 2768.9092
$ws.Range("J134").Value = 9499.5
$ws.Range("K134").Value = 8306.7276
$ws.Range("L134").Value = 28498.5
$ws.Range("M134").Value = -5771.7276
$ws.Range("N134").Value = -33568.5
$ws.Range("H136").Value = 14113.556
$ws.Range("I136").Value = 14113.556
$ws.Range("J136").Value = 0
$ws.Range("K136").Value = 42340.66800000001
$ws.Range("L136").Value = 0
$ws.Range("M136").ClearContents()
$ws.Range("N136").Value = -39790.66800000001
$ws.Range("H139").Value = 59083.816
$ws.Range("J139").Value = 59992.2
$ws.Range("L139").Value = 59992.2
$ws.Range("N139").Value = -70272.2

# ---- Sheet CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H11").Value = 4088.5
$ws.Range("J11").Value = 3161.5
$ws.Range("L11").Value = 9484.5
$ws.Range("N11").Value = -9764.5
$ws.Range("H37").Value = 83369820
$ws.Range("J37").Value = 83369820
$ws.Range("L37").Value = 250109460
$ws.Range("N37").Value = -250109684
$ws.Range("H113").Value = 2105.7693
$ws.Range("J113").Value = 2556.7646
$ws.Range("L113").Value = 7670.293799999999
$ws.Range("N113").Value = -12010.2938

# ---- Sheet GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H3").Value = 6429098
$ws.Range("I3").Value = 6667399
$ws.Range("J3").Value = 6250372
$ws.Range("K3").Value = 6667399
$ws.Range("L3").Value = 6250372
$ws.Range("M3").Value = -6667283
$ws.Range("N3").Value = -6250604
$ws.Range("H126").Value = 5842.357
$ws.Range("I126").Value = 3998.25
$ws.Range("J126").Value = 8301.166999999999
$ws.Range("K126").Value = 11994.75
$ws.Range("L126").Value = 24903.501
$ws.Range("M126").Value = -9524.75
$ws.Range("N126").Value = -29843.501

# ---- Sheet LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 13642.409
$ws.Range("I7").Value = 12487.5
$ws.Range("K7").Value = 12487.5
$ws.Range("M7").Value = -12375.5
$ws.Range("H22").Value = 1807.3
$ws.Range("I22").Value = 1744.5
$ws.Range("J22").Value = 1849.1666
$ws.Range("K22").Value = 1744.5
$ws.Range("L22").Value = 1849.1666
$ws.Range("M22").Value = -1449.5
$ws.Range("N22").Value = -2439.1666
$ws.Range("H27").Value = 1807.3
$ws.Range("I27").Value = 1744.5
$ws.Range("J27").Value = 1849.1666
$ws.Range("K27").Value = 1744.5
$ws.Range("L27").Value = 1849.1666
$ws.Range("M27").Value = -1637.5
$ws.Range("N27").Value = -2063.1666
$ws.Range("H55").Value = 221.71428
$ws.Range("I55").Value = 127.77778
$ws.Range("K55").Value = 127.77778
$ws.Range("M55").Value = 45.22221999999999
$ws.Range("H126").Value = 13642.409
$ws.Range("I126").Value = 12487.5
$ws.Range("K126").Value = 37462.5
$ws.Range("M126").Value = -34992.5
$ws.Range("H132").Value = 23319.25
$ws.Range("I132").Value = 30771.691
$ws.Range("K132").Value = 92315.073
$ws.Range("M132").Value = -89785.073

# ---- Sheet WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H123").Value = 84998.5
$ws.Range("J123").Value = 84998.5
$ws.Range("L123").Value = 84998.5
$ws.Range("N123").Value = -94798.5
$ws.Range("H132").Value = 2604.4546
$ws.Range("I132").Value = 2383.3333
$ws.Range("K132").Value = 7149.999899999999
$ws.Range("M132").Value = -4619.999899999999
